$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 639.62964
$ws.Range("I19").Value = 313.625
$ws.Range("J19").Value = 776.8946999999999
$ws.Range("K19").Value = 313.625
$ws.Range("L19").Value = 776.8946999999999
$ws.Range("M19").Value = -138.625
$ws.Range("N19").Value = -1126.8947

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 351.5
$ws.Range("J33").Value = 394.2857
$ws.Range("L33").Value = 394.2857
$ws.Range("N33").Value = -852.2857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2362.9546
$ws.Range("I40").Value = 2764.5833
$ws.Range("J40").Value = 1881
$ws.Range("K40").Value = 2764.5833
$ws.Range("L40").Value = 1881
$ws.Range("M40").Value = -2589.5833
$ws.Range("N40").Value = -2231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 8000
$ws.Range("J18").Value = 8000
$ws.Range("L18").Value = 8000
$ws.Range("N18").Value = -8644

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6546.4443
$ws.Range("I32").Value = 5478.6777
$ws.Range("J32").Value = 22296
$ws.Range("K32").Value = 5478.6777
$ws.Range("L32").Value = 22296
$ws.Range("M32").Value = -5191.6777
$ws.Range("N32").Value = -22870

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5155.2583
$ws.Range("I61").Value = 3746.4285
$ws.Range("K61").Value = 3746.4285
$ws.Range("M61").Value = -3534.4285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5003.3145
$ws.Range("I74").Value = 3064.4546
$ws.Range("J74").Value = 8284.462
$ws.Range("K74").Value = 3064.4546
$ws.Range("L74").Value = 8284.462
$ws.Range("M74").Value = -2190.4546
$ws.Range("N74").Value = -10032.462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5003.3145
$ws.Range("I77").Value = 3064.4546
$ws.Range("J77").Value = 8284.462
$ws.Range("K77").Value = 15322.273
$ws.Range("L77").Value = 41422.31
$ws.Range("M77").Value = -10954.273
$ws.Range("N77").Value = -50158.31

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8228.166999999999
$ws.Range("I132").Value = 5787.75
$ws.Range("J132").Value = 13109
$ws.Range("K132").Value = 17363.25
$ws.Range("L132").Value = 39327
$ws.Range("M132").Value = -14833.25
$ws.Range("N132").Value = -44387

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5155.2583
$ws.Range("I136").Value = 3746.4285
$ws.Range("K136").Value = 11239.2855
$ws.Range("M136").Value = -8689.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 628.5714
$ws.Range("I64").Value = 650
$ws.Range("J64").Value = 600
$ws.Range("K64").Value = 650
$ws.Range("L64").Value = 600
$ws.Range("M64").Value = -425
$ws.Range("N64").Value = -1050

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 628.5714
$ws.Range("I67").Value = 650
$ws.Range("J67").Value = 600
$ws.Range("K67").Value = 650
$ws.Range("L67").Value = 600
$ws.Range("M67").Value = 130
$ws.Range("N67").Value = -2160

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8028.278
$ws.Range("I105").Value = 5176.5713
$ws.Range("J105").Value = 12020.667
$ws.Range("K105").Value = 5176.5713
$ws.Range("L105").Value = 12020.667
$ws.Range("M105").Value = -3429.5713
$ws.Range("N105").Value = -15514.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 41945
$ws.Range("J141").Value = 41945
$ws.Range("L141").Value = 41945
$ws.Range("N141").Value = -52305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2287.3542
$ws.Range("I31").Value = 1577.8108
$ws.Range("J31").Value = 4674
$ws.Range("K31").Value = 1577.8108
$ws.Range("L31").Value = 4674
$ws.Range("M31").Value = -1282.8108
$ws.Range("N31").Value = -5264

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2287.3542
$ws.Range("I34").Value = 1577.8108
$ws.Range("J34").Value = 4674
$ws.Range("K34").Value = 1577.8108
$ws.Range("L34").Value = 4674
$ws.Range("M34").Value = -1375.8108
$ws.Range("N34").Value = -5078

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").ClearContents()
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 9713.777
$ws.Range("I122").Value = 3017.9412
$ws.Range("J122").Value = 21096.7
$ws.Range("K122").Value = 9053.8236
$ws.Range("L122").Value = 63290.10000000001
$ws.Range("M122").Value = -6603.8236
$ws.Range("N122").Value = -68190.10000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2996.0193
$ws.Range("I134").Value = 1811.2084
$ws.Range("J134").Value = 4011.5715
$ws.Range("K134").Value = 5433.6252
$ws.Range("L134").Value = 12034.7145
$ws.Range("M134").Value = -2898.6252
$ws.Range("N134").Value = -17104.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 30.545454
$ws.Range("J2").Value = 33.875
$ws.Range("L2").Value = 203.25
$ws.Range("N2").Value = -429.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 1071.6
$ws.Range("I10").Value = 1071.6
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 3214.8
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -3075.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 78.8125
$ws.Range("I38").Value = 26
$ws.Range("J38").Value = 96.416664
$ws.Range("K38").Value = 78
$ws.Range("L38").Value = 289.249992
$ws.Range("M38").Value = 269
$ws.Range("N38").Value = -983.249992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 3480
$ws.Range("J42").Value = 3980
$ws.Range("L42").Value = 11940
$ws.Range("N42").Value = -13008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1681.5652
$ws.Range("I129").Value = 1265.1
$ws.Range("J129").Value = 2001.9231
$ws.Range("K129").Value = 3795.3
$ws.Range("L129").Value = 6005.7693
$ws.Range("M129").Value = 1204.7
$ws.Range("N129").Value = -16005.7693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 52854.75
$ws.Range("J100").Value = 52854.75
$ws.Range("L100").Value = 52854.75
$ws.Range("N100").Value = -55018.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4964.5454
$ws.Range("I102").Value = 5536
$ws.Range("J102").Value = 2393
$ws.Range("K102").Value = 5536
$ws.Range("L102").Value = 2393
$ws.Range("M102").Value = -3914
$ws.Range("N102").Value = -5637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4076.1667
$ws.Range("I7").Value = 2858.2666
$ws.Range("J7").Value = 6106
$ws.Range("K7").Value = 2858.2666
$ws.Range("L7").Value = 6106
$ws.Range("M7").Value = -2746.2666
$ws.Range("N7").Value = -6330

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 6750
$ws.Range("J20").Value = 6750
$ws.Range("L20").Value = 6750
$ws.Range("N20").Value = -7202

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3477.257
$ws.Range("I40").Value = 3517.1304
$ws.Range("J40").Value = 3400.8333
$ws.Range("K40").Value = 3517.1304
$ws.Range("L40").Value = 3400.8333
$ws.Range("M40").Value = -3381.1304
$ws.Range("N40").Value = -3672.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").ClearContents()
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4076.1667
$ws.Range("I126").Value = 2858.2666
$ws.Range("J126").Value = 6106
$ws.Range("K126").Value = 8574.799800000001
$ws.Range("L126").Value = 18318
$ws.Range("M126").Value = -6104.799800000001
$ws.Range("N126").Value = -23258

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1661.9333
$ws.Range("J107").Value = 1905.75
$ws.Range("L107").Value = 5717.25
$ws.Range("N107").Value = -9557.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 985.5854
$ws.Range("I113").Value = 402.77777
$ws.Range("J113").Value = 1441.6957
$ws.Range("K113").Value = 1208.33331
$ws.Range("L113").Value = 4325.0871
$ws.Range("M113").Value = 961.66669
$ws.Range("N113").Value = -8665.087100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3940.121
$ws.Range("I136").Value = 1984.9722
$ws.Range("J136").Value = 6286.3
$ws.Range("K136").Value = 5954.9166
$ws.Range("L136").Value = 18858.9
$ws.Range("M136").Value = -3404.9166
$ws.Range("N136").Value = -23958.9
